# Slide 23, content placeholder shape: the paragraph
#     "    →（）は配列（文章）に残しておきたい"
# (a single run) gets split into three runs:
#     "    →（）は配列（文章）に" / "残して" / "おきたい"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# The original text uses U+00A0 (non-breaking space) for its leading
# indent, not plain ASCII spaces - build it explicitly so the written
# run matches byte-for-byte.
$nbsp = [string][char]0x00A0
$indent = "$nbsp$nbsp$nbsp$nbsp"

# The target paragraph starts at character 249 (1-based) and is 22
# characters long within the shape's full text range.
#   "<nbsp x4>→（）は配列（文章）に"  -> chars 249-263 (15 chars)
#   "残して"                            -> chars 264-266 (3 chars)
#   "おきたい"                          -> chars 267-270 (4 chars)
$run3 = $tr.Characters(267, 4)
$run3.Text = "おきたい"

$run2 = $tr.Characters(264, 3)
$run2.Text = "残して"

$run1 = $tr.Characters(249, 15)
$run1.Text = $indent + "→（）は配列（文章）に"

Write-Host "done splitting run"
